# Update bases das ligas (12-06-2024 23:38)
# The underlying match data (columns B:AD) for a handful of rows got
# reshuffled between rows while the row index in column A stayed fixed.
# Read every source row's data first (so nothing gets clobbered before
# it is used), then write all destinations.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# new row -> source row whose B:AD content it should receive
$srcRow130 = $ws.Range("B132:AD132").Value2
$srcRow131 = $ws.Range("B133:AD133").Value2
$srcRow132 = $ws.Range("B130:AD130").Value2
$srcRow133 = $ws.Range("B131:AD131").Value2
$srcRow134 = $ws.Range("B137:AD137").Value2
$srcRow136 = $ws.Range("B134:AD134").Value2
$srcRow137 = $ws.Range("B136:AD136").Value2
$srcRow142 = $ws.Range("B143:AD143").Value2
$srcRow143 = $ws.Range("B142:AD142").Value2

$ws.Range("B130:AD130").Value2 = $srcRow130
$ws.Range("B131:AD131").Value2 = $srcRow131
$ws.Range("B132:AD132").Value2 = $srcRow132
$ws.Range("B133:AD133").Value2 = $srcRow133
$ws.Range("B134:AD134").Value2 = $srcRow134
$ws.Range("B136:AD136").Value2 = $srcRow136
$ws.Range("B137:AD137").Value2 = $srcRow137
$ws.Range("B142:AD142").Value2 = $srcRow142
$ws.Range("B143:AD143").Value2 = $srcRow143
